$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add new observacao note about qualificacao_id = 36
$ws.Range("C3").Value = "qualificacao_id = 36 inexistente. Foi inserida como (‘36’, NULL) no banco de dados para evitar de perder essa linha"

# Row 6: add new observacao note about pais_id = 367 (and related ids), wrapped
$ws.Range("C6").Value = "pais_id = 367 inexistente. Foi inserida como (‘367’, NULL) no banco de dados para evitar de perder essa linha`nidem para pais_id = 150`nidem para pais_id = 449`nidem para pais_id = 678`nidem para pais_id = 359"
$ws.Range("C6").Font.Name = "Arial"
$ws.Range("C6").Font.Size = 10
$ws.Range("C6").WrapText = $true

# Column C needs to be wider to accommodate the longer notes
$ws.Columns(3).ColumnWidth = 90.35

# Row 6 grows taller because of the wrapped multi-line note
$ws.Rows(6).RowHeight = 58.9

# Move the active selection to C7 (also scrolls the view)
$ws.Range("C7").Select()
